$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: studies
# ---------------------------------------------------------------
$studies = $wb.Worksheets.Item("studies")
$studies.Range("B1").Value2 = "study_label"
$studies.Range("C1").Value2 = "description"
$studies.Range("D1").Value2 = "access_level"
$studies.Range("E1").Value2 = "contributors"
$studies.Range("F1").Value2 = "reference"
$studies.Range("G1").Value2 = "reference_year"

# ---------------------------------------------------------------
# Sheet: surveys
# ---------------------------------------------------------------
$surveys = $wb.Worksheets.Item("surveys")

# capture originals before any mutation
$origH1 = $surveys.Range("H1").Value2
$origI1 = $surveys.Range("I1").Value2
$origJ1 = $surveys.Range("J1").Value2
$origK1 = $surveys.Range("K1").Value2
$origJ2 = $surveys.Range("J2").Value2
$origG2 = $surveys.Range("G2").Value2

$surveys.Range("A1").Value2 = "study_id"

# apply text format before writing date-like strings so Excel does not
# reinterpret them as date serials
$surveys.Range("I1:L1").NumberFormat = "@"
$surveys.Range("K2").NumberFormat = "@"
$surveys.Range("L2").NumberFormat = "@"

# introduce the brand-new header strings in the same left-to-right order
# they first appear in the finished sheet, then backfill the carried-over
# values from right to left so source cells are read before being
# overwritten
$surveys.Range("G1").Value2 = "location_method"
$surveys.Range("H1").Style = "Normal"
$surveys.Range("H1").Value2 = "location_notes"
$surveys.Range("L1").Value2 = "time_method"

$surveys.Range("M1").Value2 = $origK1      # time_notes
$surveys.Range("K1").Value2 = $origJ1      # collection_day
$surveys.Range("J1").Value2 = $origI1      # collection_end
$surveys.Range("I1").Value2 = $origH1      # collection_start

$surveys.Range("H2").Style = "Normal"
$surveys.Range("H2").Value2 = $origG2      # "example data" -> location_notes
$surveys.Range("M2").Value2 = $origG2      # "example data" -> time_notes col
$surveys.Range("K2").Value2 = $origJ2      # "2020-01-01"
$surveys.Range("G2").Clear()
$surveys.Range("J2").Clear()
$surveys.Range("L2").Value2 = ""

# ---------------------------------------------------------------
# Sheet: counts
# ---------------------------------------------------------------
$counts = $wb.Worksheets.Item("counts")
$counts.Range("A1").Value2 = "study_id"
$counts.Range("B1").Value2 = "survey_id"
$counts.Range("C2").Value2 = "crt:72-76:CVIET"
$counts.Range("C3").Value2 = "crt:72-74:CVI"
$counts.Range("C4").Value2 = "crt:72-74:C/A_V/A_I"
$counts.Range("C5").Value2 = "crt:72-74:C|A|D_V|A|D_I"
$counts.Range("C6").Value2 = "crt:74-76:IET"

# ---------------------------------------------------------------
# Sheet: studies (continued) - do this last so that the new
# "public" string is appended to the shared-string table after
# the strings introduced by the surveys/counts edits above
# ---------------------------------------------------------------
$studies.Range("C2").ClearContents()
$studies.Range("D2").Value2 = "public"

# ---------------------------------------------------------------
# View / selection changes
# ---------------------------------------------------------------
$studies.Select()
$studies.Range("D4").Select()

$surveys.Select()
$surveys.Range("A1:M2").Select()

$counts.Select()
$counts.Range("D10").Select()

$prevalence = $wb.Worksheets.Item("prevalence")
$prevalence.Select()
$prevalence.Range("E15").Select()

$prevalence2 = $wb.Worksheets.Item("prevalence2")
$prevalence2.Select()
$prevalence2.Range("G11").Select()

$studies.Select()
